$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Mongoose session (previously a single combined row) is now split
# across two course sessions: row 7 keeps only the introduction, and row 8
# (which used to just say "Continuer les exercices") picks up the rest.
$ws.Range("C7").Value = "[Mongoose](introduction_mongoose.md)"
$ws.Range("C8").Value = "[Mongoose - la suite](mongoose2.md)<br/>[Simulacre Mongoose](simulacre_mongoose.md)"

$ws.Range("D7").Value = "[Exercice 5 - Mongoose](exercice5_mongoose.md)"
$ws.Range("D8").Value = "[Exercice 6 - Mongoose](exercice6_mongoose.md)"

# Move the active selection, matching the author's saved cursor position.
$ws.Range("C9").Select() | Out-Null
